$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4444.5
$ws.Range("I62").Value = 4444.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4444.5
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3820.5
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 4444.5
$ws.Range("I65").Value = 4444.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22222.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19102.5
$ws.Range("N65").ClearContents()
$ws.Range("H98").Value = 2805.125
$ws.Range("I98").Value = 2805.125
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2805.125
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1307.125
$ws.Range("N98").ClearContents()
$ws.Range("H112").Value = 2144.2
$ws.Range("I112").Value = 699.8570999999999
$ws.Range("J112").Value = 2450.5757
$ws.Range("K112").Value = 2099.5713
$ws.Range("L112").Value = 7351.7271
$ws.Range("M112").Value = -991.5712999999996
$ws.Range("N112").Value = -9567.7271
$ws.Range("H122").Value = 2805.125
$ws.Range("I122").Value = 2805.125
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8415.375
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5965.375
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H125").Value = 1095
$ws.Range("I125").Value = 1040
$ws.Range("J125").Value = 1150
$ws.Range("K125").Value = 9360
$ws.Range("L125").Value = 10350
$ws.Range("M125").Value = -6900
$ws.Range("N125").Value = -15270
$ws.Range("H126").Value = 29983.637
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 29983.637
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 29983.637
$ws.Range("N126").Value = -39863.637
$ws.Range("H127").Value = 782.4
$ws.Range("I127").Value = 641.875
$ws.Range("J127").Value = 1344.5
$ws.Range("K127").Value = 1925.625
$ws.Range("L127").Value = 4033.5
$ws.Range("M127").Value = 3034.375
$ws.Range("N127").Value = -13953.5
$ws.Range("H128").Value = 35600
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 35600
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 35600
$ws.Range("N128").Value = -45560
$ws.Range("H129").Value = 848.0769
$ws.Range("I129").Value = 603.5
$ws.Range("J129").Value = 892.5454999999999
$ws.Range("K129").Value = 1810.5
$ws.Range("L129").Value = 2677.6365
$ws.Range("M129").Value = 3189.5
$ws.Range("N129").Value = -12677.6365
$ws.Range("H131").Value = 1763.3334
$ws.Range("I131").Value = 2145
$ws.Range("J131").Value = 1000
$ws.Range("K131").Value = 6435
$ws.Range("L131").Value = 3000
$ws.Range("M131").Value = -1395
$ws.Range("N131").Value = -13080
$ws.Range("H137").Value = 1451.7
$ws.Range("I137").Value = 1012.63635
$ws.Range("J137").Value = 1988.3334
$ws.Range("K137").Value = 3037.90905
$ws.Range("L137").Value = 5965.0002
$ws.Range("M137").Value = -487.9090500000002
$ws.Range("N137").Value = -11065.0002

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1219.2325
$ws.Range("I31").Value = 882.75
$ws.Range("J31").Value = 1511.826
$ws.Range("K31").Value = 882.75
$ws.Range("L31").Value = 1511.826
$ws.Range("M31").Value = -587.75
$ws.Range("N31").Value = -2101.826
$ws.Range("H34").Value = 1219.2325
$ws.Range("I34").Value = 882.75
$ws.Range("J34").Value = 1511.826
$ws.Range("K34").Value = 882.75
$ws.Range("L34").Value = 1511.826
$ws.Range("M34").Value = -680.75
$ws.Range("N34").Value = -1915.826
$ws.Range("H58").Value = 892.5143
$ws.Range("I58").Value = 943.25
$ws.Range("J58").Value = 781.8182
$ws.Range("K58").Value = 943.25
$ws.Range("L58").Value = 781.8182
$ws.Range("M58").Value = -740.25
$ws.Range("N58").Value = -1187.8182
$ws.Range("H122").Value = 977.6
$ws.Range("I122").Value = 962.6667
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2888.0001
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -438.0001000000002
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 5832.5415
$ws.Range("I132").Value = 6098.45
$ws.Range("J132").Value = 4503
$ws.Range("K132").Value = 18295.35
$ws.Range("L132").Value = 13509
$ws.Range("M132").Value = -15765.35
$ws.Range("N132").Value = -18569
$ws.Range("H136").Value = 892.5143
$ws.Range("I136").Value = 943.25
$ws.Range("J136").Value = 781.8182
$ws.Range("K136").Value = 2829.75
$ws.Range("L136").Value = 2345.4546
$ws.Range("M136").Value = -279.75
$ws.Range("N136").Value = -7445.4546

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1556.9642
$ws.Range("I5").Value = 1345.0454
$ws.Range("J5").Value = 2334
$ws.Range("K5").Value = 4035.1362
$ws.Range("L5").Value = 7002
$ws.Range("M5").Value = -3923.1362
$ws.Range("N5").Value = -7226
$ws.Range("H113").Value = 563.65216
$ws.Range("I113").Value = 431.14285
$ws.Range("J113").Value = 621.625
$ws.Range("K113").Value = 1293.42855
$ws.Range("L113").Value = 1864.875
$ws.Range("M113").Value = 876.5714499999999
$ws.Range("N113").Value = -6204.875
$ws.Range("H135").Value = 1556.9642
$ws.Range("I135").Value = 1345.0454
$ws.Range("J135").Value = 2334
$ws.Range("K135").Value = 12105.4086
$ws.Range("L135").Value = 21006
$ws.Range("M135").Value = -9570.408599999999
$ws.Range("N135").Value = -26076

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1374367.1
$ws.Range("I107").Value = 2404537.5
$ws.Range("J107").Value = 806.6667
$ws.Range("K107").Value = 2404537.5
$ws.Range("L107").Value = 806.6667
$ws.Range("M107").Value = -2402617.5
$ws.Range("N107").Value = -4646.6667
$ws.Range("H122").Value = 2858.1428
$ws.Range("I122").Value = 1702.1
$ws.Range("J122").Value = 3909.0908
$ws.Range("K122").Value = 5106.299999999999
$ws.Range("L122").Value = 11727.2724
$ws.Range("M122").Value = -2656.299999999999
$ws.Range("N122").Value = -16627.2724

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("H122").Value = 20837040
$ws.Range("I122").Value = 50003104
$ws.Range("J122").Value = 4139.7144
$ws.Range("K122").Value = 150009312
$ws.Range("L122").Value = 12419.1432
$ws.Range("M122").Value = -150006862
$ws.Range("N122").Value = -17319.1432
$ws.Range("H132").Value = 79769.69500000001
$ws.Range("I132").Value = 3749.5
$ws.Range("J132").Value = 113556.445
$ws.Range("K132").Value = 11248.5
$ws.Range("L132").Value = 340669.335
$ws.Range("M132").Value = -8718.5
$ws.Range("N132").Value = -345729.335

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 59429
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 59429
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 59429
$ws.Range("N46").Value = -59891
$ws.Range("H134").Value = 59429
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 59429
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 178287
$ws.Range("N134").Value = -183357
